# The date column (A) was being stored as a real Excel date serial (with a
# custom "yyyy-mm-dd h:mm:ss"-style number format applied), but the
# downstream dataloader only consumes a plain YYYYMMDD integer. Convert every
# data row's date serial in column A to an YYYYMMDD integer and drop the
# now-unneeded custom date style so the cell goes back to the default
# (unformatted) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRows = $ws.UsedRange.Rows.Count
if ($usedRows -lt 2) { $usedRows = 2 }

for ($r = 2; $r -le $usedRows; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $serial = $cell.Value2
    if ($null -eq $serial -or $serial -eq "") { continue }

    $asDate = [DateTime]::FromOADate([double]$serial)
    $yyyymmdd = [int]$asDate.ToString("yyyyMMdd")

    $cell.Value = $yyyymmdd
    $cell.Style = "Normal"
}
